# Rescale the Customer_Reviews column (C2:C51) by a constant factor.
# The underlying data-generation process for Customer_Reviews was
# re-run with a different scaling coefficient, so every existing value
# in the column is multiplied by the same constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 0.7811488943493342

for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    $cell.Value = $v * $factor
}
